# Update the cryptos price list: refresh Price (D) and Volume(1h) (E) columns,
# plus a row swap (ShibaInu <-> WrappedliquidstakedEther2.0) for rows 15/16.
# For "Price" cells whose new text looks like a plain number (e.g. "253.57"),
# force the cell to Text format before writing so Excel keeps it as a string
# (matching the source data, which stores these as text) instead of silently
# converting it to a numeric value; ClearFormats() afterwards removes the
# now-unneeded custom number format again so no extra cell style is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '97.760.49'
$ws.Range("E2").Value = '  +3.54%  '
$ws.Range("D3").Value = '3.339.54'
$ws.Range("E3").Value = '  +8.43%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '253.57'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +7.00%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '621.97'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.79%  '
$ws.Range("E7").Value = '  +7.01%  '
$ws.Range("E8").Value = '  +1.71%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").Value = '3.335.77'
$ws.Range("E10").Value = '  +8.42%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.799'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.73%  '
$ws.Range("E12").Value = '  +1.20%  '
$ws.Range("D13").Value = '97.551.89'
$ws.Range("E13").Value = '  +3.67%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.71'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +5.19%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000245'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.92%  '
$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").Value = '3.944.34'
$ws.Range("E16").Value = '  +8.10%  '
$ws.Range("E17").Value = '  +2.31%  '
$ws.Range("D18").Value = '3.342.14'
$ws.Range("E18").Value = '  +8.07%  '
$ws.Range("E19").Value = '  +0.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.69'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '478.95'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +8.44%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.86'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.99%  '
$ws.Range("E23").Value = '  +8.89%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.07'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.86%  '
$ws.Range("E25").Value = '  +2.56%  '
$ws.Range("E26").Value = '  +3.38%  '
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("D28").Value = '3.571.15'
$ws.Range("E28").Value = '  +10.27%  '
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("E30").Value = '  +5.74%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.247'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.06%  '
$ws.Range("E32").Value = '  -0.49%  '
$ws.Range("E33").Value = '  +0.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '9.14'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '27.10'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +6.72%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '519.41'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +8.55%  '
$ws.Range("E37").Value = '  -1.65%  '
$ws.Range("E38").Value = '  -5.92%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.94'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +2.43%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '24.78'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +3.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.446'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.00%  '
$ws.Range("E42").Value = '  -0.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.68'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.788'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +17.18%  '
$ws.Range("E45").Value = '  +3.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '160.72'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.45%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.92'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +5.67%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '45.49'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +4.35%  '
$ws.Range("E50").Value = '  +6.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.49'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +5.50%  '
